# Change the table style on the "Source / Definition / Advantages" table
# from the deck's custom "Table_0" style to the built-in table style
# {628C806B-56E4-4D4F-9F87-B2BDF297116C} (previously
# {AE23D92C-B546-4E51-B0D9-429CBED5301B}).

$p = $ppt.ActivePresentation

$targetStyleId = "{628C806B-56E4-4D4F-9F87-B2BDF297116C}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
